$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three Cypher queries in column B (shared strings) ---

# Row 2 (Cases query): append an ORDER BY / LIMIT clause
$b2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

# Row 3 (Samples query): append an ORDER BY / LIMIT clause
$b3 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

# Row 4 (Files query): replace the existing trailing "order by" clause
$b4 = $ws.Range("B4").Value2
$b4 = $b4 -replace "    order by f\.file_name$", "     order By f.file_name ASC LIMIT 100"
$ws.Range("B4").Value2 = $b4

# --- Row heights grow to fit the longer wrapped query text ---
$ws.Rows.Item(2).RowHeight = 360
$ws.Rows.Item(3).RowHeight = 374.4

# --- Selection / scroll position moves to C4 ---
$ws.Range("C4").Select()
